$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Musculoskeletal Radiology (mk)) - update mean/SD/sample size
$ws.Range("B2").Value = 0.003537936797399456
$ws.Range("C2").Value = 0.02166382793468516
$ws.Range("D2").Value = 100

# Row 3 (Computed Tomography (ct)) - only sample size changes
$ws.Range("D3").Value = 96

# Row 4 - category renamed from "Emergency Radiology (er)" to "Breast Imaging (br)"
$ws.Range("A4").Value = "Breast Imaging (br)"
$ws.Range("B4").Value = 0.1432560135011233
$ws.Range("C4").Value = 0.1084180057565769
$ws.Range("D4").Value = 98

# Row 5 - category renamed from "Biomarkers and Quantative imaging (bq)" to "Geritourinary Radiology (gu)"
$ws.Range("A5").Value = "Geritourinary Radiology (gu)"
$ws.Range("B5").Value = 0.0009685406801488272
$ws.Range("C5").Value = 0.009588063559035733
$ws.Range("D5").Value = 99

# Row 6 - category renamed from "Safety and Quality (sq)" to "Ultrasound (us)"
$ws.Range("A6").Value = "Ultrasound (us)"
$ws.Range("B6").Value = 0.005598125815857243
$ws.Range("C6").Value = 0.02730579349881621
$ws.Range("D6").Value = 96

# Row 7 - category renamed from "Professionalism (pr)" to "Chest Radiology (ch)"
$ws.Range("A7").Value = "Chest Radiology (ch)"
$ws.Range("B7").Value = 0.01240836581265891
$ws.Range("C7").Value = 0.02990742360884703
$ws.Range("D7").Value = 99

# Row 8 - category renamed from "Nuclear Medicine (nm)" to "Interventional Radiology (ir)"
$ws.Range("A8").Value = "Interventional Radiology (ir)"
$ws.Range("B8").Value = 0.005393984377847852
$ws.Range("C8").Value = 0.0274330692286452
$ws.Range("D8").Value = 96

# Rows 9-11 (Chest Radiology (ch), Breast Imaging (br), Geritourinary Radiology (gu)) removed -
# the dataset now only has 7 categories instead of 10
$ws.Range("A9:D11").Delete()
